$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SignUp_Positive_01): Condition flips Yes->No, Email/Firstname/Lastname get new randomized values
$ws.Range("C2").Value = "No"
$ws.Range("I2").Value = "prabhaAutoTtMV7494@mailinator.com"
$ws.Range("K2").Value = "PrabhaAutoPjKz"
$ws.Range("L2").Value = "automationuntz"

# Row 3 (SignUp_Positive_02): Condition flips Yes->No
$ws.Range("C3").Value = "No"

# Row 5 (SignUp_Positive_04): now a chronic-disease test instead of mental-health
$ws.Range("B5").Value = "To verify the employee's selection for the 'Lunch and Learnt' chronic program."""
$ws.Range("C5").Value = "Yes"
$ws.Range("I5").Value = "prabhaAutoErsX9791@mailinator.com"
$ws.Range("K5").Value = "PrabhaAutooOBA"
$ws.Range("L5").Value = "automationgxpl"
$ws.Range("U5").Value = "Chronic Disease"

# Update the active selection to D6, matching the saved view state
$ws.Range("D6").Select()
